# Update "想去人数" (column F) figures for matching rows on both the
# "展览" and "全部类型" worksheets, as per upstream data refresh.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1375
    6  = 233
    11 = 4625
    12 = 6893
    14 = 59
    18 = 4139
    19 = 710
    22 = 2728
    26 = 367
    33 = 68
    34 = 340
    35 = 86
    40 = 160
    42 = 18
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
